$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 1968.8889
$ws.Range("I58").Value = 70.8
$ws.Range("J58").Value = 4341.5
$ws.Range("K58").Value = 212.4
$ws.Range("L58").Value = 13024.5
$ws.Range("M58").Value = -62.39999999999998
$ws.Range("N58").Value = -13324.5
$ws.Range("H64").Value = 28335.1
$ws.Range("I64").Value = 3121.4285
$ws.Range("J64").Value = 87167
$ws.Range("K64").Value = 3121.4285
$ws.Range("L64").Value = 87167
$ws.Range("M64").Value = -2873.4285
$ws.Range("N64").Value = -87663
$ws.Range("H67").Value = 28335.1
$ws.Range("I67").Value = 3121.4285
$ws.Range("J67").Value = 87167
$ws.Range("K67").Value = 3121.4285
$ws.Range("L67").Value = 87167
$ws.Range("M67").Value = -2263.4285
$ws.Range("N67").Value = -88883
$ws.Range("H116").Value = 7890.7617
$ws.Range("I116").Value = 22200.6
$ws.Range("J116").Value = 3418.9375
$ws.Range("K116").Value = 22200.6
$ws.Range("L116").Value = 3418.9375
$ws.Range("M116").Value = -18758.6
$ws.Range("N116").Value = -10302.9375
$ws.Range("H129").Value = 657.75
$ws.Range("J129").Value = 1000
$ws.Range("L129").Value = 3000
$ws.Range("N129").Value = -13000
$ws.Range("H137").Value = 33442.47
$ws.Range("I137").Value = 67741.47
$ws.Range("J137").Value = 6364.316
$ws.Range("K137").Value = 203224.41
$ws.Range("L137").Value = 19092.948
$ws.Range("M137").Value = -200674.41
$ws.Range("N137").Value = -24192.948
$ws.Range("H138").Value = 1184.7333
$ws.Range("I138").Value = 634.8302
$ws.Range("J138").Value = 1972.4324
$ws.Range("K138").Value = 1904.4906
$ws.Range("L138").Value = 5917.2972
$ws.Range("M138").Value = 3235.5094
$ws.Range("N138").Value = -16197.2972
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1777.05
$ws.Range("I61").Value = 1926.7576
$ws.Range("J61").Value = 1071.2858
$ws.Range("K61").Value = 1926.7576
$ws.Range("L61").Value = 1071.2858
$ws.Range("M61").Value = -1714.7576
$ws.Range("N61").Value = -1495.2858
$ws.Range("H63").Value = 3000
$ws.Range("I63").Value = 3000
$ws.Range("K63").Value = 3000
$ws.Range("M63").Value = -2314
$ws.Range("H66").Value = 3000
$ws.Range("I66").Value = 3000
$ws.Range("K66").Value = 15000
$ws.Range("M66").Value = -11568
$ws.Range("H80").Value = 16463.076
$ws.Range("J80").Value = 17418.25
$ws.Range("L80").Value = 17418.25
$ws.Range("N80").Value = -19414.25
$ws.Range("H83").Value = 16463.076
$ws.Range("J83").Value = 17418.25
$ws.Range("L83").Value = 52254.75
$ws.Range("N83").Value = -62238.75
$ws.Range("H132").Value = 2737823
$ws.Range("I132").Value = 3191904.8
$ws.Range("J132").Value = 1123309.4
$ws.Range("K132").Value = 9575714.399999999
$ws.Range("L132").Value = 3369928.2
$ws.Range("M132").Value = -9573184.399999999
$ws.Range("N132").Value = -3374988.2
$ws.Range("H136").Value = 1777.05
$ws.Range("I136").Value = 1926.7576
$ws.Range("J136").Value = 1071.2858
$ws.Range("K136").Value = 5780.2728
$ws.Range("L136").Value = 3213.8574
$ws.Range("M136").Value = -3230.2728
$ws.Range("N136").Value = -8313.857400000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 5595.4287
$ws.Range("I82").Value = 3480.8333
$ws.Range("J82").Value = 18283
$ws.Range("K82").Value = 3480.8333
$ws.Range("L82").Value = 18283
$ws.Range("M82").Value = -3097.8333
$ws.Range("N82").Value = -19049
$ws.Range("H85").Value = 5595.4287
$ws.Range("I85").Value = 3480.8333
$ws.Range("J85").Value = 18283
$ws.Range("K85").Value = 3480.8333
$ws.Range("L85").Value = 18283
$ws.Range("M85").Value = -2154.8333
$ws.Range("N85").Value = -20935
$ws.Range("H134").Value = 21647.703
$ws.Range("I134").Value = 1050.1111
$ws.Range("J134").Value = 62842.89
$ws.Range("K134").Value = 3150.3333
$ws.Range("L134").Value = 188528.67
$ws.Range("M134").Value = -615.3333000000002
$ws.Range("N134").Value = -193598.67
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7553.95
$ws.Range("I31").Value = 5631.5884
$ws.Range("J31").Value = 18447.334
$ws.Range("K31").Value = 5631.5884
$ws.Range("L31").Value = 18447.334
$ws.Range("M31").Value = -5336.5884
$ws.Range("N31").Value = -19037.334
$ws.Range("H34").Value = 7553.95
$ws.Range("I34").Value = 5631.5884
$ws.Range("J34").Value = 18447.334
$ws.Range("K34").Value = 5631.5884
$ws.Range("L34").Value = 18447.334
$ws.Range("M34").Value = -5429.5884
$ws.Range("N34").Value = -18851.334
$ws.Range("H41").Value = 8742.857
$ws.Range("I41").Value = 7750
$ws.Range("J41").Value = 14700
$ws.Range("K41").Value = 7750
$ws.Range("L41").Value = 14700
$ws.Range("M41").Value = -7322
$ws.Range("N41").Value = -15556
$ws.Range("H50").Value = 10245
$ws.Range("J50").Value = 10245
$ws.Range("L50").Value = 10245
$ws.Range("N50").Value = -11495
$ws.Range("H51").Value = 9562
$ws.Range("J51").Value = 9562
$ws.Range("L51").Value = 9562
$ws.Range("N51").Value = -11034
$ws.Range("H58").Value = 2458.3845
$ws.Range("I58").Value = 756.63416
$ws.Range("J58").Value = 8801.272000000001
$ws.Range("K58").Value = 756.63416
$ws.Range("L58").Value = 8801.272000000001
$ws.Range("M58").Value = -553.63416
$ws.Range("N58").Value = -9207.272000000001
$ws.Range("H59").Value = 13955.444
$ws.Range("J59").Value = 13955.444
$ws.Range("L59").Value = 13955.444
$ws.Range("N59").Value = -16245.444
$ws.Range("H60").Value = 9501.25
$ws.Range("J60").Value = 9501.25
$ws.Range("L60").Value = 9501.25
$ws.Range("N60").Value = -10523.25
$ws.Range("H61").Value = 9562
$ws.Range("J61").Value = 9562
$ws.Range("L61").Value = 9562
$ws.Range("N61").Value = -10258
$ws.Range("H62").Value = 2132.5
$ws.Range("I62").Value = 1742
$ws.Range("J62").Value = 2783.3333
$ws.Range("K62").Value = 1742
$ws.Range("L62").Value = 2783.3333
$ws.Range("M62").Value = -1118
$ws.Range("N62").Value = -4031.3333
$ws.Range("H65").Value = 2132.5
$ws.Range("I65").Value = 1742
$ws.Range("J65").Value = 2783.3333
$ws.Range("K65").Value = 8710
$ws.Range("L65").Value = 13916.6665
$ws.Range("M65").Value = -5590
$ws.Range("N65").Value = -20156.6665
$ws.Range("H68").Value = 13973.8
$ws.Range("I68").Value = 268
$ws.Range("J68").Value = 17400.25
$ws.Range("K68").Value = 268
$ws.Range("L68").Value = 17400.25
$ws.Range("M68").Value = 481
$ws.Range("N68").Value = -18898.25
$ws.Range("H71").Value = 13973.8
$ws.Range("I71").Value = 268
$ws.Range("J71").Value = 17400.25
$ws.Range("K71").Value = 804
$ws.Range("L71").Value = 52200.75
$ws.Range("M71").Value = 2940
$ws.Range("N71").Value = -59688.75
$ws.Range("H74").Value = 13065.5
$ws.Range("I74").Value = 285
$ws.Range("J74").Value = 14227.363
$ws.Range("K74").Value = 285
$ws.Range("L74").Value = 14227.363
$ws.Range("M74").Value = 589
$ws.Range("N74").Value = -15975.363
$ws.Range("H77").Value = 13065.5
$ws.Range("I77").Value = 285
$ws.Range("J77").Value = 14227.363
$ws.Range("K77").Value = 855
$ws.Range("L77").Value = 42682.089
$ws.Range("M77").Value = 3513
$ws.Range("N77").Value = -51418.089
$ws.Range("H134").Value = 1525.2632
$ws.Range("I134").Value = 1392
$ws.Range("K134").Value = 4176
$ws.Range("M134").Value = -1641
$ws.Range("H136").Value = 2458.3845
$ws.Range("I136").Value = 756.63416
$ws.Range("J136").Value = 8801.272000000001
$ws.Range("K136").Value = 2269.90248
$ws.Range("L136").Value = 26403.816
$ws.Range("M136").Value = 280.0975200000003
$ws.Range("N136").Value = -31503.816
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 420
$ws.Range("I122").Value = 350
$ws.Range("J122").Value = 431.05264
$ws.Range("K122").Value = 3150
$ws.Range("L122").Value = 3879.47376
$ws.Range("M122").Value = -700
$ws.Range("N122").Value = -8779.473760000001
$ws.Range("H131").Value = 34723070
$ws.Range("I131").Value = 450
$ws.Range("J131").Value = 40323492
$ws.Range("K131").Value = 1350
$ws.Range("L131").Value = 120970476
$ws.Range("M131").Value = 3690
$ws.Range("N131").Value = -120980556
$ws.Range("H137").Value = 8429392
$ws.Range("I137").Value = 15152234
$ws.Range("J137").Value = 5067970
$ws.Range("K137").Value = 45456702
$ws.Range("L137").Value = 15203910
$ws.Range("M137").Value = -45451602
$ws.Range("N137").Value = -15214110
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 28000
$ws.Range("J39").Value = 28000
$ws.Range("L39").Value = 28000
$ws.Range("N39").Value = -29064
$ws.Range("H107").Value = 19936.4
$ws.Range("I107").Value = 420.5
$ws.Range("J107").Value = 98000
$ws.Range("K107").Value = 420.5
$ws.Range("L107").Value = 98000
$ws.Range("M107").Value = 1499.5
$ws.Range("N107").Value = -101840
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2043.4736
$ws.Range("I68").Value = 1783.3334
$ws.Range("J68").Value = 2489.4285
$ws.Range("K68").Value = 1783.3334
$ws.Range("L68").Value = 2489.4285
$ws.Range("M68").Value = -1034.3334
$ws.Range("N68").Value = -3987.4285
$ws.Range("H71").Value = 2043.4736
$ws.Range("I71").Value = 1783.3334
$ws.Range("J71").Value = 2489.4285
$ws.Range("K71").Value = 8916.666999999999
$ws.Range("L71").Value = 12447.1425
$ws.Range("M71").Value = -5172.666999999999
$ws.Range("N71").Value = -19935.1425
$ws.Range("H93").Value = 2918.1428
$ws.Range("I93").Value = 3003
$ws.Range("J93").Value = 2904
$ws.Range("K93").Value = 3003
$ws.Range("L93").Value = 2904
$ws.Range("M93").Value = -1755
$ws.Range("N93").Value = -5400
$ws.Range("H100").Value = 45873.652
$ws.Range("I100").Value = 251751
$ws.Range("J100").Value = 2531.0527
$ws.Range("K100").Value = 251751
$ws.Range("L100").Value = 2531.0527
$ws.Range("M100").Value = -251210
$ws.Range("N100").Value = -3613.0527
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1293417.9
$ws.Range("I136").Value = 1429066.2
$ws.Range("J136").Value = 771693.4
$ws.Range("K136").Value = 4287198.6
$ws.Range("L136").Value = 2315080.2
$ws.Range("M136").Value = -4284648.6
$ws.Range("N136").Value = -2320180.2
